# Atualização automática via cronjob
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column A sequence numbers for existing rows (2-7 shift by the
#     insertion of the new record below; row 8 keeps its original 0) ---
$ws.Cells.Item(3, 1).Value = 3
$ws.Cells.Item(4, 1).Value = 4
$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(6, 1).Value = 6
$ws.Cells.Item(7, 1).Value = 7

# --- Update estoque_atualizado (column G) values that changed ---
$ws.Cells.Item(4, 7).Value = 483
$ws.Cells.Item(5, 7).Value = 200
$ws.Cells.Item(8, 7).Value = 119

# --- Append new row 9 with the new sale record ---
# "Dia" (B) and "id_produto" (E) must stay plain text like the rest of the
# column (values such as "2025-04-17" / "000032" are identifiers, not real
# dates/numbers), so force text format before assigning, then drop back to
# the default style so the cell matches the look of its siblings.
$ws.Cells.Item(9, 2).NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = "2025-04-17"
$ws.Cells.Item(9, 2).Style = "Normal"

$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "000032"
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(9, 3).Value = 84
$ws.Cells.Item(9, 4).Value = "METALURGICA SATO DA AMAZONIA LTDA"
$ws.Cells.Item(9, 6).Value = "LIMPADOR VEJA MULTIUSO GOLD 500ML"
$ws.Cells.Item(9, 7).Value = 826
$ws.Cells.Item(9, 8).Value = $false

# match the bordered/bold/centered style used by the other data rows in
# column A (font bold, thin border all sides, centered, top-aligned)
$ws.Cells.Item(9, 1).Font.Bold = $true
$ws.Cells.Item(9, 1).HorizontalAlignment = -4108
$ws.Cells.Item(9, 1).VerticalAlignment = -4160
$ws.Cells.Item(9, 1).Borders.LineStyle = 1
$ws.Cells.Item(9, 1).Value = 2
